$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.792.81'
$ws.Range('E2').Value = '  -2.74%  '
$ws.Range('D3').Value = '1.615.72'
$ws.Range('E3').Value = '  -3.13%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.003'
$ws.Range('E5').Value = '  +0.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '306.57'
$ws.Range('E6').Value = '  -1.96%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3907'
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3803'
$ws.Range('E8').Value = '  -3.23%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.002'
$ws.Range('E9').Value = '  -0.14%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.350'
$ws.Range('E10').Value = '  -3.35%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '48.46'
$ws.Range('E11').Value = '  -6.13%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08430'
$ws.Range('E12').Value = '  -2.16%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '23.73'
$ws.Range('E13').Value = '  -5.84%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.975'
$ws.Range('E14').Value = '  -4.06%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001270'
$ws.Range('E15').Value = '  -3.51%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.400'
$ws.Range('E16').Value = '  -3.80%  '
$ws.Range('D17').Value = '1.610.48'
$ws.Range('E17').Value = '  -3.53%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '93.32'
$ws.Range('E18').Value = '  +0.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06909'
$ws.Range('E19').Value = '  -1.99%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.88'
$ws.Range('E20').Value = '  -4.32%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.787'
$ws.Range('E21').Value = '  -3.68%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.003'
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.36'
$ws.Range('E23').Value = '  -4.03%  '
$ws.Range('D24').Value = '23.825.25'
$ws.Range('E24').Value = '  -2.60%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.413'
$ws.Range('E25').Value = '  +1.51%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.827'
$ws.Range('E26').Value = '  +3.46%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.06'
$ws.Range('E27').Value = '  -4.65%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '157.04'
$ws.Range('E28').Value = '  -2.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '138.88'
$ws.Range('E29').Value = '  -5.92%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.226'
$ws.Range('E30').Value = '  -10.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.732'
$ws.Range('E31').Value = '  -6.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.474'
$ws.Range('E32').Value = '  -1.70%  '
$ws.Range('D33').Value = '1.795.90'
$ws.Range('E33').Value = '  -3.31%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08011'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9532'
$ws.Range('E35').Value = '  -2.47%  '
$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02852'
$ws.Range('E36').Value = '  -5.22%  '
$ws.Range('B37').Value = 'InternetComputer(DFINITY)'
$ws.Range('C37').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.541'
$ws.Range('E37').Value = '  -6.10%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2637'
$ws.Range('E38').Value = '  -5.85%  '
$ws.Range('E39').Value = '  -3.17%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '10.38'
$ws.Range('E40').Value = '  +1.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '13.31'
$ws.Range('E41').Value = '  -1.35%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.418'
$ws.Range('E42').Value = '  -7.16%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.7420'
$ws.Range('E43').Value = '  -5.80%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '15.78'
$ws.Range('E44').Value = '  -3.55%  '
$ws.Range('E45').Value = '  -4.10%  '
$ws.Range('E46').Value = '  -4.46%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.050'
$ws.Range('E47').Value = '  -2.78%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.003'
$ws.Range('E48').Value = '  +0.03%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.08210'
$ws.Range('E49').Value = '  -4.15%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '132.02'
$ws.Range('E50').Value = '  -3.82%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.216'
$ws.Range('E51').Value = '  -7.65%  '
